$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Huambo_City"
$ws.Range("A3").Value = "Luanda_City"

$ws.Range("C4").Select()
